$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = 'Última actualización: 18:35:28'
$ws1.Cells.Item(3,1).Value = 'Total filas: 406'

$sheet1Data = @(
    @(195, '12:35:30', '13:02', '14_ABASTO', 27, 'LP1912'),
    @(196, '12:01:11', '13:02', '15_ABASTO', 61, 'LP1912'),
    @(208, '11:46:46', '13:26', '15_ABASTO', 100, 'LP1912'),
    @(209, '11:46:46', '13:26', '14_ABASTO', 100, 'LP1912'),
    @(220, '12:01:11', '13:51', '215A_EL PATO', 110, 'LP1912'),
    @(221, '13:51:32', '13:51', '11_ETCHEVERRY', 0, 'LP1912'),
    @(283, '14:59:23', '16:09', '27_EL RETIRO', 70, 'LP1912'),
    @(284, '14:49:07', '16:09', '14_ABASTO', 80, 'LP1912'),
    @(295, '16:34:19', '16:34', '23_HERNANDEZ', 0, 'LP1912'),
    @(296, '16:34:19', '16:34', '11_ETCHEVERRY', 0, 'LP1912'),
    @(323, '15:36:13', '17:21', '26_HERNANDEZ', 105, 'LP1912'),
    @(324, '16:20:15', '17:21', '16_SANTA ANA', 61, 'LP1912'),
    @(325, '17:12:54', '17:21', '10_OLMOS', 9, 'LP1912'),
    @(335, '16:45:34', '17:38', '27_EL RETIRO', 53, 'LP1912'),
    @(336, '15:59:02', '17:38', '17_ROMERO', 99, 'LP1912'),
    @(348, '15:59:02', '17:52', '81_EL PELIGRO', 113, 'LP1912'),
    @(349, '17:36:40', '17:52', '11_ETCHEVERRY', 16, 'LP1912'),
    @(360, '18:13:12', '18:21', '16_SANTA ANA', 8, 'LP1912'),
    @(361, '16:34:19', '18:21', '26_HERNANDEZ', 107, 'LP1912'),
    @(372, '18:35:28', '18:36', '23_HERNANDEZ', 1, 'LP1912'),
    @(374, '17:50:30', '18:40', '15_ABASTO', 50, 'LP1912'),
    @(375, '18:13:12', '18:41', '14_ABASTO', 28, 'LP1912'),
    @(376, '17:12:54', '18:42', '14_ABASTO', 90, 'LP1912'),
    @(377, '16:53:02', '18:45', '14_ABASTO', 112, 'LP1912'),
    @(378, '17:12:54', '18:47', '14X44_ABASTO', 95, 'LP1912'),
    @(379, '16:53:02', '18:48', '14X44_ABASTO', 115, 'LP1912'),
    @(380, '18:13:12', '18:52', '15_ABASTO', 39, 'LP1912'),
    @(381, '17:50:30', '18:56', '10_OLMOS', 66, 'LP1912'),
    @(382, '17:12:54', '18:58', '215A_EL PATO', 106, 'LP1912'),
    @(383, '18:13:12', '18:59', '215A_EL PATO', 46, 'LP1912'),
    @(384, '17:12:54', '19:04', '11_ETCHEVERRY', 112, 'LP1912'),
    @(385, '17:36:40', '19:04', '23_HERNANDEZ', 88, 'LP1912'),
    @(387, '18:13:12', '19:05', '11_ETCHEVERRY', 52, 'LP1912'),
    @(388, '17:12:54', '19:10', '16_P MOR-SANTA ANA', 118, 'LP1912'),
    @(389, '18:13:12', '19:11', '16_P MOR-SANTA ANA', 58, 'LP1912'),
    @(390, '18:35:28', '19:12', '10_OLMOS', 37, 'LP1912'),
    @(391, '17:36:40', '19:16', '27_EL RETIRO', 100, 'LP1912'),
    @(392, '18:35:28', '19:16', '15_ABASTO', 41, 'LP1912'),
    @(393, '18:13:12', '19:17', '27_EL RETIRO', 64, 'LP1912'),
    @(394, '17:50:30', '19:20', '14_ABASTO', 90, 'LP1912'),
    @(395, '18:35:28', '19:20', '16_SANTA ANA', 45, 'LP1912'),
    @(396, '17:36:40', '19:21', '26_HERNANDEZ', 105, 'LP1912'),
    @(397, '18:13:12', '19:28', '15_ABASTO', 75, 'LP1912'),
    @(398, '17:36:40', '19:30', '225_GOMEZ', 114, 'LP1912'),
    @(399, '18:35:28', '19:30', '16_SANTA ANA', 55, 'LP1912'),
    @(400, '17:50:30', '19:40', '215C_EL PATO', 110, 'LP1912'),
    @(401, '18:13:12', '19:50', '11X44_ETCHEVERRY', 97, 'LP1912'),
    @(402, '18:35:28', '19:50', '16_P MOR-SANTA ANA', 75, 'LP1912'),
    @(403, '18:13:12', '19:51', '81_EL PELIGRO', 98, 'LP1912'),
    @(404, '18:13:12', '19:51', '16_P MOR-SANTA ANA', 98, 'LP1912'),
    @(405, '18:13:12', '19:59', '17_ROMERO', 106, 'LP1912'),
    @(406, '18:35:28', '20:10', '16_P MOR-167 Y 521', 95, 'LP1912'),
    @(407, '18:13:12', '20:11', '16_P MOR-167 Y 521', 118, 'LP1912'),
    @(408, '18:35:28', '20:21', '26_HERNANDEZ', 106, 'LP1912'),
    @(409, '18:35:28', '20:22', '11_ETCHEVERRY', 107, 'LP1912'),
    @(410, '18:35:28', '20:23', '215A_EL PATO', 108, 'LP1912'),
    @(411, '18:35:28', '20:31', '225_GOMEZ', 116, 'LP1912')
)

foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r,1).Value = $row[1]
    $ws1.Cells.Item($r,2).Value = $row[2]
    $ws1.Cells.Item($r,3).Value = $row[3]
    $ws1.Cells.Item($r,4).Value = $row[4]
    $ws1.Cells.Item($r,5).Value = $row[5]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = 'Última actualización: 18:35:28'
$ws2.Cells.Item(3,1).Value = 'Total filas: 42'

$sheet2Data = @(
    @(47, '18:35:28', '20:23', '215A_EL PATO', 108, 'LP1912')
)

foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r,1).Value = $row[1]
    $ws2.Cells.Item($r,2).Value = $row[2]
    $ws2.Cells.Item($r,3).Value = $row[3]
    $ws2.Cells.Item($r,4).Value = $row[4]
    $ws2.Cells.Item($r,5).Value = $row[5]
}

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = 'Última actualización: 18:35:28'

Write-Output "Done"